$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Remove the per-file category labels that used to live in column D
# (dns.pcap, netflow.pcap, sip_rtp.pcap, SkypeIRC.cap). The chart only
# ever referenced columns E and F, so these labels are unused and are
# being dropped as part of the "removed unused diagrams" cleanup.
$ws.Range("D5:D8").ClearContents()

# Update the active selection left behind in the sheet view.
$ws.Range("O9").Select()
